$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the standings table (rows 5-18). Two new games were reported for
#    several teams, so win/loss/points totals and the team ordering by rank
#    changed. Rank numbers in column B stay 1..14 (unaffected).
# ---------------------------------------------------------------------------

function Set-StandingsRow {
    param($RowNum, $Team, $Games, $Wins, $Losses, $Score, $Points)
    $ws.Cells.Item($RowNum, 3).Value2 = $Team     # C - Команда
    $ws.Cells.Item($RowNum, 4).Value2 = $Games    # D - Игры
    $ws.Cells.Item($RowNum, 5).Value2 = $Wins     # E - Побед
    $ws.Cells.Item($RowNum, 6).Value2 = $Losses   # F - Поражений
    $ws.Cells.Item($RowNum, 7).Value2 = $Score    # G - Мячи
    $ws.Cells.Item($RowNum, 8).Value2 = $Points   # H - Очки
}

Set-StandingsRow 5  "ISsoft"                  12 10 2  "818 - 660" 22
Set-StandingsRow 6  "Эра-Недвижимости плюс"   12 10 2  "892 - 720" 22
Set-StandingsRow 7  "БГУФК"                   12 9  3  "835 - 668" 21
Set-StandingsRow 8  "Грушвиль"                12 9  3  "973 - 788" 21
Set-StandingsRow 9  "ОПЛАТИ"                  12 9  3  "901 - 726" 21
Set-StandingsRow 10 "GOLDEN HILL"             12 8  4  "851 - 804" 20
Set-StandingsRow 11 "Mapogo males"            12 7  5  "896 - 873" 19
Set-StandingsRow 12 "SIRIUS"                  12 6  6  "798 - 695" 18
Set-StandingsRow 13 "Стрела"                  12 4  8  "723 - 772" 16
Set-StandingsRow 14 "VSS"                     12 4  8  "727 - 837" 16
Set-StandingsRow 15 "Eagles"                  12 3  9  "686 - 755" 15
Set-StandingsRow 16 "NORD"                    12 3  9  "646 - 946" 15
Set-StandingsRow 17 "ЛФК"                     12 2  10 "676 - 831" 14
Set-StandingsRow 18 "Минск 7х"                12 0  12 "524 - 871" 12

Write-Host "standings updated"

# ---------------------------------------------------------------------------
# 2. Append the two new match days (15 and 16 March 2025) below the existing
#    results list, matching the established layout: a merged date-header row
#    (style copied from an existing date row) followed by one merged row per
#    match result (style copied from an existing match row).
# ---------------------------------------------------------------------------

function Add-DateHeaderRow {
    param($RowNum, $SerialDate, $TemplateRow)
    $target = "B" + $RowNum + ":H" + $RowNum
    $ws.Range($target).Merge()
    $ws.Range("B" + $TemplateRow + ":H" + $TemplateRow).Copy()
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
    $ws.Range("B" + $RowNum).Value2 = $SerialDate
}

function Add-MatchRow {
    param($RowNum, $Text, $TemplateRow)
    $target = "B" + $RowNum + ":H" + $RowNum
    $ws.Range($target).Merge()
    $ws.Range("B" + $TemplateRow + ":H" + $TemplateRow).Copy()
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
    $ws.Range("B" + $RowNum).Value2 = $Text
    $ws.Rows.Item($RowNum).RowHeight = 19.95
}

# 15 March 2025 (serial 45731)
Add-DateHeaderRow 120 45731 115
Add-MatchRow 121 "БГУФК - Mapogo males 79:65 (16:30, БНТУ)" 116
Add-MatchRow 122 "Минск 7х - ЛФК 37:63 (18:00, БНТУ)" 116
Add-MatchRow 123 "Грушвиль - Eagles 68:50 (19:30, БНТУ)" 116

# 16 March 2025 (serial 45732)
Add-DateHeaderRow 124 45732 115
Add-MatchRow 125 "ISsoft - NORD 93:74 (11:00, БНТУ)" 116
Add-MatchRow 126 "Стрела - SIRIUS 52:67 (12:30, БНТУ)" 116
Add-MatchRow 127 "Эра-Недвижимости плюс - GOLDEN HILL 74:65 (14:00, БНТУ)" 116
Add-MatchRow 128 "VSS - ОПЛАТИ 38:91 (15:30, БНТУ)" 116

Write-Host "new match rows added"
